# airJetLab.xlsx edit: formulas up to task 4 corrected (ABS() removed from the
# velocity formula so negative pressure-head readings correctly surface as
# #NUM! errors instead of being silently rectified) and the R/L/H/V column
# headers in row 14 retyped (now carrying/losing trailing spaces per column,
# matching how the author retyped them by hand), plus plot/view tidy-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 headers: retype the R / L / H / V labels for each of the three
#     x-station blocks (A:D, F:I, K:N). The new text no longer carries the
#     "(mm)" suffix, and trailing spaces now vary column to column exactly
#     as the author left them.
$ws.Range("A14").Value2 = "R "
$ws.Range("B14").Value2 = "L "
$ws.Range("C14").Value2 = "H"
$ws.Range("D14").Value2 = "V"

$ws.Range("F14").Value2 = "R "
$ws.Range("G14").Value2 = "L "
$ws.Range("H14").Value2 = "H "
$ws.Range("I14").Value2 = "V"

$ws.Range("K14").Value2 = "R"
$ws.Range("L14").Value2 = "L "
$ws.Range("M14").Value2 = "H "
$ws.Range("N14").Value2 = "V"

# --- Velocity formulas: drop the ABS() that was masking negative
#     pressure-head readings. Rows 15 hold standalone formulas; rows
#     16-29/16-35/16-39 are shared-formula blocks, so each block's master
#     cell is re-entered across the whole range to keep the sharing intact.
$ws.Range("D15").Formula = "=40*SQRT(SIN(0.2251474735)*10*C15/1000)"
$ws.Range("I15").Formula = "=40*SQRT(SIN(0.2251474735)*10*H15/1000)"
$ws.Range("N15").Formula = "=40*SQRT(SIN(0.2251474735)*10*M15/1000)"

$ws.Range("D16:D29").Formula = "=40*SQRT(SIN(0.2251474735)*10*C16/1000)"
$ws.Range("I16:I35").Formula = "=40*SQRT(SIN(0.2251474735)*10*H16/1000)"
$ws.Range("N16:N39").Formula = "=40*SQRT(SIN(0.2251474735)*10*M16/1000)"

# --- View tidy-up: scroll/select so F19 is the active cell, matching the
#     author's final on-screen state.
$ws.Range("F19").Select()

$wb.Save()
